{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the four content changes described in the diff:\n//  1. Split/rewrite the \"peripheral IV\" paragraphs into three paragraphs.\n//  2. Insert three new paragraphs describing the port chamber before\n//     \"Surgically placed under the skin\".\n//  3. Extend the suture sentence with an extra clause.\n//  4. Change \"two days\" to \"the first week\" in the post-op lifting guidance.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nfunction findIndexByText(targetText) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.trim() === targetText) {\n      return i;\n    }\n  }\n  return -1;\n}\n\n// --- Change 1 -------------------------------------------------------\n// Old paragraph 1: \"A Peripheral intravenous line is placed before each\n// dose or doses of chemotherapy or immunotherapy and removed that same day\"\n// Old paragraph 2: \"A peripheral IV is not suitable for FLOT chemotherapy,\n// which requires a 24-hour infusion\"\nconst idx1 = findIndexByText(\n  \"A Peripheral intravenous line is placed before each dose or doses of chemotherapy or immunotherapy and removed that same day\"\n);\nif (idx1 === -1) {\n  throw new Error(\"Could not find the peripheral IV intro paragraph\");\n}\nconst pLine = items[idx1];\npLine.insertText(\n  \"A Peripheral intravenous line IV is placed through the skin into a vein in the hand or forearm before each dose or doses of chemotherapy or immunotherapy.\",\n  \"Replace\"\n);\n// New paragraph inserted right after it.\nconst pCatheter = pLine.insertParagraph(\n  \"The catheter is removed the same day and then replaced when it\\u2019s time for the next dose, which is typically one to three weeks later.\",\n  \"After\"\n);\n\nconst idx2 = findIndexByText(\n  \"A peripheral IV is not suitable for FLOT chemotherapy, which requires a 24-hour infusion\"\n);\nif (idx2 === -1) {\n  throw new Error(\"Could not find the FLOT chemotherapy paragraph\");\n}\nitems[idx2].insertText(\n  \"For patients who receive FLOT chemotherapy, a peripheral IV is not suitable because this treatment requires a 24-hour infusion of drug.\",\n  \"Replace\"\n);\n\n// --- Change 2 -------------------------------------------------------\n// Insert three new paragraphs right before \"Surgically placed under the skin\".\nconst idx3 = findIndexByText(\"Surgically placed under the skin\");\nif (idx3 === -1) {\n  throw new Error(\"Could not find the 'Surgically placed under the skin' paragraph\");\n}\nconst pSurgically = items[idx3];\nconst pChamber = pSurgically.insertParagraph(\n  \"It contains a small chamber with a flexible rubber top and a tube that goes into the veins near the heart\",\n  \"Before\"\n);\nconst pNeedle = pChamber.insertParagraph(\n  \"When it comes time to administer drugs, a needle is passed through the skin into the chamber\",\n  \"After\"\n);\npNeedle.insertParagraph(\n  \"This avoids having to find a vein underneath the skin in the arm or hand. The port is ideal for those with small veins which might be difficult to access with a peripheral intravenous line\",\n  \"After\"\n);\n\n// --- Change 3 -------------------------------------------------------\n// Extend the suture sentence.\nconst idx4 = findIndexByText(\"The skin is closed with sutures that dissolve on their own\");\nif (idx4 === -1) {\n  throw new Error(\"Could not find the sutures paragraph\");\n}\nitems[idx4].insertText(\n  \"The skin is closed with sutures that dissolve on their own and don\\u2019t need to be removed.\",\n  \"Replace\"\n);\n\n// --- Change 4 -------------------------------------------------------\n// \"two days\" -> \"the first week\"\nconst idx5 = findIndexByText(\n  \"We recommend no lifting for two days to avoid bruising in the area\"\n);\nif (idx5 === -1) {\n  throw new Error(\"Could not find the no-lifting paragraph\");\n}\nitems[idx5].insertText(\n  \"We recommend no lifting for the first week to avoid bruising in the area\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the four content changes described in the diff:\n#  1. Split/rewrite the \"peripheral IV\" paragraphs into three paragraphs.\n#  2. Insert three new paragraphs describing the port chamber before\n#     \"Surgically placed under the skin\".\n#  3. Extend the suture sentence with an extra clause.\n#  4. Change \"two days\" to \"the first week\" in the post-op lifting guidance.\n\n$d = $word.ActiveDocument\n$cr = [char]13\n$rsquo = [char]0x2019\n\nfunction Get-ParagraphIndexByText($doc, $text) {\n    $i = 0\n    foreach ($p in $doc.Paragraphs) {\n        $i = $i + 1\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# --- Change 1 ---------------------------------------------------------\n# Old paragraph 1: \"A Peripheral intravenous line is placed before each\n# dose or doses of chemotherapy or immunotherapy and removed that same day\"\n$idx1 = Get-ParagraphIndexByText $d \"A Peripheral intravenous line is placed before each dose or doses of chemotherapy or immunotherapy and removed that same day\"\nif ($idx1 -eq -1) {\n    throw \"Could not find the peripheral IV intro paragraph\"\n}\n$p1 = $d.Paragraphs.Item($idx1)\n$p1.Range.Text = \"A Peripheral intravenous line IV is placed through the skin into a vein in the hand or forearm before each dose or doses of chemotherapy or immunotherapy.\"\n\n# Insert the new \"catheter\" paragraph right after it.\n$p1again = $d.Paragraphs.Item($idx1)\n$p1again.Range.InsertAfter(\"The catheter is removed the same day and then replaced when it\" + $rsquo + \"s time for the next dose, which is typically one to three weeks later.\" + $cr)\n\n# Old paragraph 2: \"A peripheral IV is not suitable for FLOT chemotherapy,\n# which requires a 24-hour infusion\"\n$idx2 = Get-ParagraphIndexByText $d \"A peripheral IV is not suitable for FLOT chemotherapy, which requires a 24-hour infusion\"\nif ($idx2 -eq -1) {\n    throw \"Could not find the FLOT chemotherapy paragraph\"\n}\n$p2 = $d.Paragraphs.Item($idx2)\n$p2.Range.Text = \"For patients who receive FLOT chemotherapy, a peripheral IV is not suitable because this treatment requires a 24-hour infusion of drug.\"\n\n# --- Change 2 -----------------------------------------------------------\n# Insert three new paragraphs right before \"Surgically placed under the skin\".\n$idx3 = Get-ParagraphIndexByText $d \"Surgically placed under the skin\"\nif ($idx3 -eq -1) {\n    throw \"Could not find the 'Surgically placed under the skin' paragraph\"\n}\n$p3 = $d.Paragraphs.Item($idx3)\n$newParasText = \"It contains a small chamber with a flexible rubber top and a tube that goes into the veins near the heart\" + $cr + \"When it comes time to administer drugs, a needle is passed through the skin into the chamber\" + $cr + \"This avoids having to find a vein underneath the skin in the arm or hand. The port is ideal for those with small veins which might be difficult to access with a peripheral intravenous line\" + $cr\n$p3.Range.InsertBefore($newParasText)\n\n# --- Change 3 -------------------------------------------------------------\n# Extend the suture sentence.\n$idx4 = Get-ParagraphIndexByText $d \"The skin is closed with sutures that dissolve on their own\"\nif ($idx4 -eq -1) {\n    throw \"Could not find the sutures paragraph\"\n}\n$p4 = $d.Paragraphs.Item($idx4)\n$p4.Range.Text = \"The skin is closed with sutures that dissolve on their own and don\" + $rsquo + \"t need to be removed.\"\n\n# --- Change 4 ---------------------------------------------------------------\n# \"two days\" -> \"the first week\"\n$idx5 = Get-ParagraphIndexByText $d \"We recommend no lifting for two days to avoid bruising in the area\"\nif ($idx5 -eq -1) {\n    throw \"Could not find the no-lifting paragraph\"\n}\n$p5 = $d.Paragraphs.Item($idx5)\n$p5.Range.Text = \"We recommend no lifting for the first week to avoid bruising in the area\"\n"}
